$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = 0.6458333333333334
$ws.Range("D17").Value = 0.7209302325581395

$ws.Range("B18").Value = 0.8444444444444444
$ws.Range("C18").Value = 0.6909090909090909
$ws.Range("D18").Value = 0.7599999999999999

$ws.Range("B19").Value = 0.7419354838709677
$ws.Range("C19").Value = 0.7419354838709677
$ws.Range("D19").Value = 0.7419354838709677
$ws.Range("E19").Value = 0.7419354838709677

$ws.Range("B20").Value = 0.7451388888888889
$ws.Range("C20").Value = 0.7533492822966508
$ws.Range("D20").Value = 0.7404651162790696

$ws.Range("B21").Value = 0.7632915173237754
$ws.Range("C21").Value = 0.7419354838709677
$ws.Range("D21").Value = 0.7440360090022505

$ws.Range("B22").Value = 0.6521739130434783
$ws.Range("D22").Value = 0.7142857142857143

$ws.Range("B23").Value = 0.8297872340425532
$ws.Range("C23").Value = 0.7090909090909091
$ws.Range("D23").Value = 0.764705882352941

$ws.Range("B24").Value = 0.7419354838709677
$ws.Range("C24").Value = 0.7419354838709677
$ws.Range("D24").Value = 0.7419354838709677
$ws.Range("E24").Value = 0.7419354838709677

$ws.Range("B25").Value = 0.7409805735430157
$ws.Range("C25").Value = 0.7492822966507178
$ws.Range("D25").Value = 0.7394957983193277

$ws.Range("B26").Value = 0.7572140491182
$ws.Range("C26").Value = 0.7419354838709677
$ws.Range("D26").Value = 0.7441040932502031
